# Apply the "Add files via upload" edit to JenkinTrigger.xlsx:
#  - Bold + box-border the header row (A1:C1)
#  - Box-border the data rows (A2:C22)
#  - Widen column C
#  - Flip a few YES/NO values in column C (C2:C5 -> YES, C6/C14/C16 -> NO)
#  - Move the active selection to A6 and scroll the view back to the top
#  - Force the print page orientation to Portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting -----------------------------------------------------------
# Build the bold+border combo on a scratch cell first and paste the format
# down onto the header row in a single combined write, then apply the
# plain border to the rest of the table.
$scratch = $ws.Range("E1")
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1
$scratch.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$scratch.Clear()
$excel.CutCopyMode = $false

$ws.Range("A2:C22").Borders.LineStyle = 1

$ws.Columns.Item(3).ColumnWidth = 10.2857

# --- Data updates -----------------------------------------------------------
$ws.Range("C2").Value = "YES"
$ws.Range("C3").Value = "YES"
$ws.Range("C4").Value = "YES"
$ws.Range("C5").Value = "YES"
$ws.Range("C6").Value = "NO"
$ws.Range("C14").Value = "NO"
$ws.Range("C16").Value = "NO"

# --- View state -------------------------------------------------------------
$null = $ws.Range("A6").Select()

# --- Page setup ---------------------------------------------------------
$ws.PageSetup.Orientation = 1
